$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 28
$srcRow = 27

# Copy formatting (style) from the previous data row so the new row matches
$ws.Range("A$srcRow`:H$srcRow").Copy() | Out-Null
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Cells.Item($newRow, 1).Value = "2025-08-18 09:44:43 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-18 15:14:43 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
